# Text updates as supplied by PM&C.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# Shorten the "Short title" text.
$ws.Range("B2").Value = "Increase Indigenous enrolment in early childhood education"

# Split the old "Source: ..." cell into a "Source" label (A10) plus the
# bare citation text (B10), matching the label/value layout used by the
# other rows on this sheet (e.g. A8 "Notes" / B8 ...).
$ws.Range("A10").Value = "Source"
$ws.Range("B10").Value = "ABS unpublished, Preschool Education, Australia 2015; Australian Demographic Statistics, June 2015; Births, Australia, 2015; Estimates and Projections, Aboriginal and Torres Strait Islander Australians, 2001 to 2026; Schools, Australia 2015"
